$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, pushing existing rows 49-54 down to 50-55.
$ws.Rows.Item(49).Insert()

# Fill in the new row 49 with the new record. Columns A,B,C,E,F,G,H,I,N,Q,R
# carry the same constant values as the surrounding rows in this block.
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 45106
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112026
$ws.Cells.Item(49, 7).Value = "Haba"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 100
$ws.Cells.Item(49, 11).Value = 16000
$ws.Cells.Item(49, 12).Value = 17000
$ws.Cells.Item(49, 13).Value = 16500
$ws.Cells.Item(49, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 16).Value = 660
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
